$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.164.50'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '2.652.13'
$ws.Range('E3').Value = '  +3.36%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '592.45'
$ws.Range('E5').Value = '  +2.37%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.94'
$ws.Range('E6').Value = '  -0.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.586'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '2.650.74'
$ws.Range('E9').Value = '  +3.30%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.66'
$ws.Range('E11').Value = '  +2.51%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.152'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.354'
$ws.Range('E13').Value = '  +1.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.33'
$ws.Range('E14').Value = '  +2.06%  '
$ws.Range('D15').Value = '3.131.77'
$ws.Range('E15').Value = '  +3.45%  '
$ws.Range('D16').Value = '63.119.39'
$ws.Range('E16').Value = '  +0.86%  '
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').Value = '2.649.69'
$ws.Range('E18').Value = '  +4.09%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.40'
$ws.Range('E19').Value = '  +3.03%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.39'
$ws.Range('E20').Value = '  +1.53%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '338.19'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.77'
$ws.Range('E22').Value = '  +1.92%  '
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.43'
$ws.Range('E24').Value = '  +0.90%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.68'
$ws.Range('E25').Value = '  +6.19%  '
$ws.Range('E26').Value = '  +1.52%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.52'
$ws.Range('E27').Value = '  +0.96%  '
$ws.Range('B28').Value = 'Binance-PegBSC-USD'
$ws.Range('C28').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  +0.21%  '
$ws.Range('B29').Value = 'InternetComputer(DFINITY)'
$ws.Range('C29').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.41'
$ws.Range('E29').Value = '  +2.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '535.86'
$ws.Range('E30').Value = '  +18.25%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.80'
$ws.Range('E31').Value = '  -1.27%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.83'
$ws.Range('E32').Value = '  +12.87%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.98'
$ws.Range('E33').Value = '  +3.16%  '
$ws.Range('D34').Value = '0.0₃0806'
$ws.Range('E34').Value = '  +1.71%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '173.06'
$ws.Range('E35').Value = '  -1.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.04'
$ws.Range('E36').Value = '  +13.94%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.406'
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('E38').Value = '  -0.06%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.03'
$ws.Range('E39').Value = '  +1.49%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.82'
$ws.Range('E40').Value = '  +8.42%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '172.43'
$ws.Range('E41').Value = '  +9.66%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').Value = '  +0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.12'
$ws.Range('E43').Value = '  -0.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.75'
$ws.Range('E44').Value = '  +2.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '22.02'
$ws.Range('E45').Value = '  +5.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0558'
$ws.Range('E46').Value = '  +4.97%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.632'
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('B48').Value = 'Stellar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0960'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0239'
$ws.Range('E49').Value = '  +2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.67'
$ws.Range('E50').Value = '  +4.34%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.36'
$ws.Range('E51').Value = '  -0.47%  '
